$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("MAIN_CONTROLLER")
$ws2 = $wb.Worksheets.Item("DATASHEET")
$ws3 = $wb.Worksheets.Item("MOBILE_CONFIGURATION")

# ----- DATASHEET: introduce the new test-data rows first (FOS_PostSanction /
# FOS3UW_to_postSanction.xlsx / CPC_1stTouchPoint_Approval /
# CPC_1stTouchPoint_Approval.xlsx / FOS8.xlsx), in this order -----
$ws2.Range("A5").Value = 2
$ws2.Range("B5").Value = "Y"
$ws2.Range("C5").Value = "FOS_PostSanction"
$ws2.Range("D5").Value = "FOS3UW_to_postSanction.xlsx"
$ws2.Range("E5").Value = 2
$ws2.Range("F5").Value = 20

$ws2.Range("A6").Value = 2
$ws2.Range("B6").Value = "Y"
$ws2.Range("C6").Value = "CPC_1stTouchPoint_Approval"
$ws2.Range("C6").Font.Color = 0
$ws2.Range("D6").Value = "CPC_1stTouchPoint_Approval.xlsx"
$ws2.Range("E6").Value = 20
$ws2.Range("F6").Value = 20

$ws2.Range("D3").Value = "FOS8.xlsx"
$ws2.Range("F4").Value = 20

# ----- MAIN_CONTROLLER -----
$ws1.Range("B2").Value = "N"
$ws1.Range("B4").Value = "N"

# Row5 used to carry the red-fill highlight style; it is no longer highlighted.
$ws1.Range("B5").ClearFormats() | Out-Null
$ws1.Range("D5").Value = "FOS_PostSanction"
$ws1.Range("E5").Value = "FOS_PostSanction"

# New row for the CPC 1st touch point approval process.
$ws1.Range("A6").Value = 3
$ws1.Range("B6").Value = "Y"
$ws1.Range("C6").Value = "local"
$ws1.Range("D6").Value = "CPC_1stTouchPoint_Approval"
$ws1.Range("E6").Value = "CPC_1stTouchPoint_Approval"

# ----- MOBILE_CONFIGURATION: duplicate row 3 into row 4 for the new
# FOS_PostSanction process -----
$ws3.Range("A4").Value = 1
$ws3.Range("B4").Value = $ws3.Range("B3").Value2
$ws3.Range("C4").Value = "FOS_PostSanction"
$ws3.Range("D4").Value = $ws3.Range("D3").Value2
$ws3.Range("E4").Value = $ws3.Range("E3").Value2
$ws3.Range("F4").Value = $ws3.Range("F3").Value2
$ws3.Range("G4").Value = $ws3.Range("G3").Value2
$ws3.Range("H4").Value = $ws3.Range("H3").Value2
$ws3.Range("I4").Value = $ws3.Range("I3").Value2
$ws3.Range("J4").Value = $ws3.Range("J3").Value2
$ws3.Range("K4").Value = $ws3.Range("K3").Value2
$ws3.Range("L4").Value = $ws3.Range("L3").Value2
$ws3.Range("M4").Value = $ws3.Range("M3").Value2
$ws3.Range("N4").Value = $ws3.Range("N3").Value2
$ws3.Range("O4").Value = $ws3.Range("O3").Value2

# Match formatting used on the row this was copied from.
$ws3.Range("D4").Style = $ws3.Range("D3").Style
$ws3.Range("E4").Style = $ws3.Range("E3").Style
$ws3.Range("H4").Style = $ws3.Range("H3").Style
$ws3.Range("K4").Style = $ws3.Range("K3").Style
$ws3.Range("N4").Style = $ws3.Range("N3").Style
$ws3.Range("O4").Style = $ws3.Range("O3").Style
$ws3.Range("P4").Style = $ws3.Range("P3").Style

# ----- Selections / active sheet -----
$ws2.Range("D13").Select() | Out-Null
$ws3.Range("E14").Select() | Out-Null
$ws1.Activate() | Out-Null
